$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.212.38"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "2.365.96"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'504.39"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'129.69"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.542"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").Value = "2.375.43"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'0.0986"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +7.08%  "
$ws.Range("D13").Value = "'0.324"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "2.789.96"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "56.206.99"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "'21.75"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "2.391.61"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "'10.00"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("D20").Value = "'310.04"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'4.02"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'65.69"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'0.370"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").Value = "'173.02"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "0.0₃0712"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "'5.84"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").Value = "'17.54"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").Value = "'1.18"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "'3.67"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").Value = "'0.827"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'36.39"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "'127.33"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").Value = "'0.562"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").Value = "'238.04"
$ws.Range("E47").Value = "  -5.16%  "
$ws.Range("D48").Value = "'0.0481"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "'0.0206"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").Value = "'16.84"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  +0.16%  "
